# Generate Report for Handoff
# Updates status text "In Translation" -> "Ready for handoff" and refreshes
# the handoff timestamps on the Overview, zh-cn, and de-de sheets, and widens
# the "Status" columns to fit the new (longer) text.
#
# Note: ColumnWidth snaps to the host's internal character-width grid
# (1/12 of a character), so 16.33 is chosen as the input that lands on the
# closest representable width to the target ~17.22 "characters".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-22 03:07:58"
$overview.Range("E:E").ColumnWidth = 16.33
$overview.Range("F:F").ColumnWidth = 16.33

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-22 03:07:53"
$zhcn.Range("C:C").ColumnWidth = 16.33

# --- de-de sheet ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-22 03:07:58"
$dede.Range("C:C").ColumnWidth = 16.33
